# Finish changes from roy's review
# Fill in the previously-placeholder (999) start_time / play_duration
# values for block_num 2 and 3 in the first run's header block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = 34.082500400000001
$ws.Range("D3").Value = 1.6276437000000001

$ws.Range("C4").Value = 49.004559700000001
$ws.Range("D4").Value = 2.6037429000000003
